$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: bump the commission percentage (H2) which ripples into F2/I2 ---
$ws.Range("H2").Value = 0.4
$ws.Range("K2").Value = 46053

# --- STT (A column) renumbering for rows 5..30 ---
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10
$ws.Range("A12").Value = 11
$ws.Range("A13").Value = 12
$ws.Range("A14").Value = 13
$ws.Range("A15").Value = 14
$ws.Range("A16").Value = 15
$ws.Range("A17").Value = 16
$ws.Range("A18").Value = 17
$ws.Range("A19").Value = 18
$ws.Range("A20").Value = 19
$ws.Range("A21").Value = 20
$ws.Range("A22").Value = 21
$ws.Range("A23").Value = 22
$ws.Range("A24").Value = 23
$ws.Range("A25").Value = 24
$ws.Range("A26").Value = 25
$ws.Range("A27").Value = 26
$ws.Range("A28").Value = 27
$ws.Range("A29").Value = 28
$ws.Range("A30").Value = 29

# --- Row 7: debt for "Nguyen Huu Nhan" got paid off in full ---
$ws.Range("E7").Formula = "=92000+8000"
$ws.Range("K7").Value = 46048
$ws.Range("M7").Value = "Đã trả đủ"
$ws.Range("B7:M7").Interior.Color = 5296274

# --- Row 8: same person, second debt line, also paid off in full ---
$ws.Range("E8").Value = 20000
$ws.Range("K8").Value = 46048
$ws.Range("M8").Value = "Đã trả đủ"
$ws.Range("B8:M8").Interior.Color = 5296274

# --- Row 9: partial payment recorded ---
$ws.Range("E9").Formula = "=22000"

# --- View state: move the active selection ---
$ws.Activate()
$ws.Range("F22").Select()
